# tree of life power point fig added
#
# Moves four existing shapes upward and adds nine new rounded-rectangle
# callout boxes (with connecting text) that describe the Filter 1 / Filter 2
# cis/trans split and the two multiple-testing correction thresholds.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# PowerPoint COM exposes Left/Top/Width/Height in points while the OOXML
# stores EMU (914400 EMU/in, 12700 EMU/pt). The host round-trips assigned
# point values through a single-precision float before truncating to EMU,
# which can shave a unit off the intended value; nudging by half an EMU
# compensates for that truncation so the saved EMU matches exactly.
function EMU([double]$v) {
    return ($v / 12700.0) + (0.5 / 12700.0)
}

function SetGeom($shape, [double]$x, [double]$y, [double]$cx, [double]$cy) {
    $shape.Left = EMU $x
    $shape.Top = EMU $y
    $shape.Width = EMU $cx
    $shape.Height = EMU $cy
}

# ---------------------------------------------------------------------
# 1. Reposition the existing flow-chart shapes.
# ---------------------------------------------------------------------

# "501 original pairs" rounded rectangle
$rr501 = $s.Shapes.Item(1)
$rr501.Top = EMU 36823

# Straight arrow connectors below it
$connLeft = $s.Shapes.Item(2)
$connLeft.Top = EMU 875373

$connRight = $s.Shapes.Item(3)
$connRight.Top = EMU 875373

# "Filter 1 (n=406)" rounded rectangle
$filter1 = $s.Shapes.Item(4)
$filter1.Left = EMU 1440873
$filter1.Top = EMU 1591804

# "Filter 2 (n=95)" rounded rectangle
$filter2 = $s.Shapes.Item(5)
$filter2.Top = EMU 1591804

# ---------------------------------------------------------------------
# 2. Add the new rounded-rectangle callouts.
#
# New shapes are built by duplicating an existing styled rounded
# rectangle (so they inherit the theme-based <p:style> block, bodyPr and
# run formatting exactly) and then overwriting geometry/name/text.
#
# PowerPoint's shape-id allocator is a monotonically increasing counter
# that is untouched by deletions, so to land on the exact ids used by the
# authored deck (11,12,13 then 17..22, skipping 14-16) a handful of
# disposable shapes are created and immediately deleted to "consume" the
# ids that must be skipped.
# ---------------------------------------------------------------------

# $sz is the OOXML <a:rPr sz="..."> value (hundredths of a point); the
# COM Font.Size property wants plain points, so divide by 100 before
# assigning. The duplicated template is already sz=1400 (14pt), so the
# common case needs no extra call.
function NewCallout([string]$name, [double]$x, [double]$y, [double]$cx, [double]$cy, [int]$sz) {
    $shape = $filter1.Duplicate().Item(1)
    SetGeom $shape $x $y $cx $cy
    $shape.Name = $name
    if ($sz -ne 1400) {
        $shape.TextFrame.TextRange.Font.Size = $sz / 100.0
    }
    return $shape
}

# -- consume ids 2,3,5,7 (already-used ids 1,4,6,8,9,10 are skipped
#    automatically by the allocator) so the first real new shape lands on id 11
for ($k = 0; $k -lt 4; $k++) {
    $tmp = $s.Shapes.AddShape(5, 0, 0, 1, 1)
    $tmp.Delete()
}

# id 11 - "Cis-cis* (n=27)"
$sh11 = NewCallout "Rounded Rectangle 10" 293056 2332660 887341 628403 1400
$sh11.TextFrame.TextRange.Text = "Cis-cis* (n=27)"

# id 12 - "Cis-trans (n=368)"
$sh12 = NewCallout "Rounded Rectangle 11" 1465303 2332661 993177 628403 1400
$sh12.TextFrame.TextRange.Text = "Cis-trans (n=368)"

# id 13 - "Trans-trans (n=11)"
$sh13 = NewCallout "Rounded Rectangle 12" 2686391 2332661 969773 628403 1200
$sh13.TextFrame.TextRange.Text = "Trans-trans (n=11)"

# -- consume ids 14,15,16 so the next real new shape lands on id 17
for ($k = 0; $k -lt 3; $k++) {
    $tmp = $s.Shapes.AddShape(5, 0, 0, 1, 1)
    $tmp.Delete()
}

# id 17 - "Cis-cis* (n=21)"
$sh17 = NewCallout "Rounded Rectangle 16" 5227323 2330397 887341 628403 1400
$sh17.TextFrame.TextRange.Text = "Cis-cis* (n=21)"

# id 18 - "Cis-trans (n=72)"
$sh18 = NewCallout "Rounded Rectangle 17" 6399570 2330398 993177 628403 1400
$sh18.TextFrame.TextRange.Text = "Cis-trans (n=72)"

# id 19 - "Trans-trans (n=2)"
$sh19 = NewCallout "Rounded Rectangle 18" 7620658 2330398 969773 628403 1200
$sh19.TextFrame.TextRange.Text = "Trans-trans (n=2)"

# id 20 - "Permutation p < 4.4e-6"
$sh20 = NewCallout "Rounded Rectangle 19" 3718177 3317229 1345228 869000 1400
$sh20.TextFrame.TextRange.Text = "Permutation p < 4.4e-6"

# id 21 - "GWAS GC Lambda F correction p < 4.48e-6"
$sh21 = NewCallout "Rounded Rectangle 20" 3718177 4338629 1345228 869000 1400
$sh21.TextFrame.TextRange.Text = "GWAS GC Lambda F correction p < 4.48e-6"

# id 22 - "GWAS GC Lambda Chisq correction p < 4.48e-6" (three runs)
$sh22 = NewCallout "Rounded Rectangle 21" 3718177 5401806 1345228 869000 1400
$tr22 = $sh22.TextFrame.TextRange
$tr22.Text = "GWAS GC Lambda "
$tr22b = $tr22.InsertAfter("Chisq")
$tr22c = $tr22b.InsertAfter(" correction p < 4.48e-6")
